$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 259 (existing last row) with corrected values
$ws.Range("C259").Value = 6243090940000
$ws.Range("D259").Value = 6243090940000
$ws.Range("E259").Value = 6243090940000
$ws.Range("F259").Value = 6243090940000

# Add three new rows of data: 260, 261, 262
$newRows = @(
    @{ Row = 260; A = 45108.41666666666; C = 6355692770000 },
    @{ Row = 261; A = 45139.41666666666; C = 6337051350000 },
    @{ Row = 262; A = 45170.41666666666; C = 6359425540000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "ECONOMICS:CZM2"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.C
    $ws.Cells.Item($row, 5).Value = $r.C
    $ws.Cells.Item($row, 6).Value = $r.C
    $ws.Cells.Item($row, 7).Value = 0
}

# Copy formatting from row 259 to the new rows (260:262), matching original style
$ws.Range("A259").Copy() | Out-Null
$ws.Range("A260:A262").PasteSpecial(-4122) | Out-Null
